$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# New test case row, appended right after the last existing row (37).
$row = 38

$ws.Cells.Item($row, 1).Value = "ENW035"
$ws.Cells.Item($row, 2).Value = "OPQA-1713"
$ws.Cells.Item($row, 3).Value = "Verify that the error should be displayed when there is not enough space in enw to accept the Neon records after clicking send to ENDNote button"
$ws.Cells.Item($row, 4).Value = "Y"

# Match formatting of the row immediately above it (borders/wrap/style).
$ws.Range("A37:E37").Copy()
$ws.Range("A38:E38").PasteSpecial(-4122) # xlPasteFormats
$ws.Rows.Item($row).RowHeight = 30

# Re-point the view the same way Excel left it after the edit.
$ws.Application.CutCopyMode = $false
$ws.Range("C38").Select()
$ws.Application.ActiveWindow.ScrollRow = 34
